$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 10 first (Delossantos.38 / S13 / C1) so row indices above it stay stable,
# then delete row 5 (Mcmurtrie.820 / M13 / H14).
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(5).Delete()

# Mirror the resulting selection captured in the saved file (entire-row selection at row 5).
$ws.Range("A5:XFD5").Select() | Out-Null
